# suppression durée traitement
# Removes the "fr-posologie.dureeTraitement" element family (and its two
# children) from the Elements sheet, bumps the Metadata "Date" value, flips
# "fr-posologie.dateDePrise" Max from 1 to *, and refreshes the
# "fr-posologie.momentDePrise" Short/Definition wording.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet: bump the IG "Date" property (row 8, column B) ---
$meta.Range("B8").Value() = "2025-07-21T07:18:37+00:00"

# --- Elements sheet ---
# Rows 37-39 hold fr-posologie.dureeTraitement / .valeur / .unite.
# Deleting them shifts the following rows (momentDePrise family) up.
$elements.Rows("37:39").Delete()

# fr-posologie.dateDePrise (now still row 36): Max 1 -> *
$elements.Range("G36").Value() = "*"
$elements.Range("AH36").Value() = "*"

# fr-posologie.momentDePrise (now row 37): reword the Short/Definition text
$newDefinition = "Définition du moment de prise au cours de la journée (ex : 30 minutes avant le repas)"
$elements.Range("L37").Value() = $newDefinition
$elements.Range("M37").Value() = $newDefinition
